$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New candidature row (2018-07-04, DEVELOPPEUR-INTEGRATEUR DE LOGICIELS @ Villeneuve-d'Ascq,
# found via Keljob / Pole Emploi).
# Copy A6's date style down to A7 first so the new date cell reuses the existing
# numFmtId="14" style instead of creating a new number format/style entry.
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = 43285

# Set text cells in the same left-to-right-ish order as the source edit so new
# shared-string entries are appended in the expected order.
$ws.Range("C7").Value = "DEVELOPPEUR-INTEGRATEUR DE LOGICIELS"
$ws.Range("F7").Value = "Villeneuve-d'Ascq"
$ws.Range("G7").Value = "Pôle Emploi via Keljob"
$ws.Range("B7").Value = "?"
$ws.Range("D7").Value = "CDI"
$ws.Range("E7").Value = "Temps Plein"

# Column width adjustments (characters, COM ColumnWidth units). The values
# below are chosen so the engine's internal pixel-quantized column width
# (round(chars*6+5) px) lands on the pixel count closest to the OOXML
# <col width="..."> targets (12.7109375 / 44 / 29.140625 "Excel chars").
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334
$ws.Columns.Item(6).ColumnWidth = 43.166666666666664
$ws.Columns.Item(7).ColumnWidth = 28.333333333333332

# Move the active selection to H6 (as in the edited workbook).
$null = $ws.Range("H6").Select()
